# Serbian language output fixes
# - Translate header row to Serbian.
# - Refresh scan_time values and switch the "present" flag column from a
#   numeric 1/0 to a textual True/False.
# - Row 13's box number was corrected and its scan_time refreshed.
# - A brand-new row 14 was scanned and appended (present = False).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row translations (row 1) ----
$ws.Range("A1").Value = "Klijent"
$ws.Range("B1").Value = "Kutija"
$ws.Range("C1").Value = "Tura"
$ws.Range("D1").Value = "Skenirao"
$ws.Range("E1").Value = "Vreme skeniranja"
$ws.Range("F1").Value = "Nalazila se u bazi"

# ---- Updated scan_time + present (True) for rows 2-12 ----
$scanTimes = @{
    2  = "2025-04-09 14:08:54"
    3  = "2025-04-09 14:08:58"
    4  = "2025-04-09 14:08:56"
    5  = "2025-04-09 14:09:00"
    6  = "2025-04-09 14:09:32"
    7  = "2025-04-09 14:09:01"
    8  = "2025-04-09 14:09:37"
    9  = "2025-04-09 14:09:49"
    10 = "2025-04-09 14:09:35"
    11 = "2025-04-09 14:09:39"
    12 = "2025-04-09 14:09:50"
}

foreach ($r in $scanTimes.Keys) {
    $ws.Cells.Item($r, 5).Value = $scanTimes[$r]
    # Leading apostrophe forces Excel to store this as text ("True"),
    # matching the workbook's textual present/absent flag instead of a
    # boolean TRUE.
    $ws.Cells.Item($r, 6).Value = "'True"
}

# ---- Row 13: box number correction, refreshed scan_time, present = False ----
$ws.Range("B13").Value = 26002680643
$ws.Range("E13").Value = "2025-04-09 14:09:06"
$ws.Range("F13").Value = "'False"

# ---- New row 14 ----
$ws.Range("A14").Value = ""
$ws.Range("B14").Value = 26002680641
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = "Milica1"
$ws.Range("E14").Value = "2025-04-09 14:09:53"
$ws.Range("F14").Value = "'False"
